$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) cells to Text format so numeric-looking strings
# such as "1.480" or "13.30" are preserved exactly as text, matching the
# original inlineStr cell content (trailing zeros / multi-dot values would
# otherwise be reinterpreted as numbers by Excel).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.265.14'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.860.07'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7104'
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '237.86'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9993'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07949'
$ws.Range('E8').Value = '  +2.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3039'
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.56'
$ws.Range('E10').Value = '  +0.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08186'
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.846.42'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.185'
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7049'
$ws.Range('E14').Value = '  -3.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.76'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.238.11'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.855'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007886'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.30'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '238.09'
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9974'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.075.78'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.415'
$ws.Range('E24').Value = '  -2.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.57'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.951'
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1439'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.12'
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.935'
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.435'
$ws.Range('E30').Value = '  +2.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.480'
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.384'
$ws.Range('E32').Value = '  -2.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.019'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05222'
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('E35').Value = '  -2.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7124'
$ws.Range('E36').Value = '  +1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9999'
$ws.Range('E37').Value = '  -2.98%  '
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01855'
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.725'
$ws.Range('E40').Value = '  +1.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9309'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.129.02'
$ws.Range('E42').Value = '  +4.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4280'
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.848'
$ws.Range('E44').Value = '  -3.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.07'
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9986'
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.97'
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5351'
$ws.Range('E48').Value = '  -4.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.768'
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.978.74'
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.169'
$ws.Range('E51').Value = '  -0.52%  '
